# Regenerate sval data to filter save games: update B2:E8 and G2:G8
# with new computed values (G = B + C + D + E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1169995834814548,  0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.104883657715537)
    3 = @(3.272327238179451,   1.626987699542094,  0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4 = @(0.1169995834814548,  0.04103571897497393,0.1496068669990043, 0.5333859586016987, 0.8410281280571317)
    5 = @(1.445647641019636,   1.626987699542094,  0.7210945179870265, 0.5333859586016987, 4.327115817150455)
    6 = @(3.272327238179451,   1.626987699542094,  0.1496068669990043, 13.86384647080068,  18.91276827552123)
    7 = @(0.6545652718822623,  1.626987699542094,  189.6080260415259,  0.5333859586016987, 192.422964971552)
    8 = @(1.445647641019636,   1.626987699542094,  3.223369029078222,  0.5333859586016987, 6.82939032824165)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E: IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G: sum
}
